$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    ,@('ussurv1062', 0, 'Manufacturing, PMI, Prices')
    ,@('ussurv1050', 0, 'Services, NMI/PMI, Prices, SA')
    ,@('ussurv1046', 0, 'NMI/PMI, New Orders, SA')
    ,@('ussurv1049', 0, 'NMI/PMI, Inventories')
    ,@('ussurv1058', 0, 'Manufacturing, PMI, Employment')
    ,@('ussurv1047', 0, 'Services, NMI/PMI, Employment')
    ,@('uslama4977', 2, 'United States, Labor Market Indicators, Labor Market Condition Index, Kansas City Fed LMCI, Level of Activity')
    ,@('uslama4978', 0, 'United States, Labor Market Indicators, Labor Market Condition Index, Kansas City Fed LMCI, Momentum')
    ,@('uslead0010', 1, 'United States, Leading Indicators, Conference Board, Business Cycle Indicators, Composite Indexes-Leading Economic Indicators, Composite Index of 10 Leading Indicators, SA, Index')
    ,@('usmost0036', 1, 'United States, Monetary Statistics, Monetary Aggregates, M2, Total, SA, USD')
    ,@('ussurv1134', 0, 'United States, Business Surveys, NFIB, Small Business Economic Trends, Small Business Optimism Index, Database, SA, Index')
    ,@('ussurv1395', 0, 'United States, Business Surveys, NFIB, Small Business Economic Trends, Expected Credit Conditions, Next Three Months, Net, Database, SA')
    ,@('ussurv1416', 2, 'Single Most Important Problem, Inflation, Database')
    ,@('ussurv1419', 2, 'Single Most Important Problem, Cost of Labour, Database')
    ,@('ussurv1417', 2, 'Single Most Important Problem, Poor Sales, Database')
    ,@('ussurv1418', 2, 'Single Most Important Problem, Financial & Interest Rates')
    ,@('ussurv1415', 2, 'Single Most Important Problem, Taxes')
    ,@('usrate0827', 0, 'TIPS, 10 Year')
    ,@('usrate0851', 0, 'TIPS, 30 Year')
    ,@('usrate0803', 0, 'TIPS, 5 Year')
    ,@('ustips5f5', 0, 'TIPS, Yield')
    ,@('usbkeven5f5', 0, 'Break-Even Inflation Rate')
    ,@('uspric0041', 1, 'United States, Import Prices, Locality of Origin, All Commodities, China, Index')
    ,@('uslama3349', 2, 'Unemployment, National, Jobless Claims, Continuing, Total')
    ,@('usgpfi0221', 1, 'Federal Government Budget, Current Receipts, Tax, on Production & Imports, Customs Duties, AR, USD')
    ,@('ussurv1397', 2, 'United States, Business Surveys, NFIB, Small Business Economic Trends, Actual Interest Rate Paid On Short-Term Loans by Borrowers, Average Interest Rate Paid, Report, SA')
)

$startRow = 28
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
# Conditional formatting: built-in "Highlight Duplicate Values" rule over
# the series column (A2:A53), matching Excel's standard red-on-red style
# (font FF9C0006 / fill FFFFC7CE). FormatCondition.Color is a VBA OLE_COLOR
# (0xBBGGRR), so the RGB bytes are reordered: 9C0006 -> 0x06009C = 393372,
# FFC7CE -> 0xCEC7FF = 13551615.
$cfRange = $ws.Range("A2:A53")
$cfRange.FormatConditions.Delete()
$fc = $cfRange.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1          # xlDuplicate (highlight duplicates, not uniques)
$fc.Font.Color = 393372     # 0x06009C -> RGB(9C,00,06)
$fc.Interior.Color = 13551615  # 0xCEC7FF -> RGB(FF,C7,CE)
$fc.Priority = 2

# View tweaks captured in the diff (zoom + active selection)
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("B17").Select()
